$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.001.42"
$ws.Range("E2").Value = "  -1.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.908.97"
$ws.Range("E3").Value = "  -3.21%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  -1.38%  "
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07712"
$ws.Range("E9").Value = "  -2.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9798"
$ws.Range("E10").Value = "  -0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.08"
$ws.Range("E11").Value = "  -2.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.890.36"
$ws.Range("E12").Value = "  -4.27%  "
$ws.Range("E13").Value = "  -2.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.932"
$ws.Range("E14").Value = "  -3.33%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07022"
$ws.Range("E15").Value = "  -1.29%  "
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -4.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009451"
$ws.Range("E18").Value = "  -4.78%  "
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.957.91"
$ws.Range("E21").Value = "  -1.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.318"
$ws.Range("E22").Value = "  -3.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.88"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.091"
$ws.Range("E24").Value = "  -0.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.12"
$ws.Range("E25").Value = "  -0.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.06"
$ws.Range("E26").Value = "  -2.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.662"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.46"
$ws.Range("E28").Value = "  -1.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.850"
$ws.Range("E29").Value = "  -1.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09287"
$ws.Range("E30").Value = "  -1.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8641"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.068"
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.248"
$ws.Range("E33").Value = "  -5.48%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.025"
$ws.Range("E34").Value = "  -3.38%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05741"
$ws.Range("E35").Value = "  -1.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.156"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02042"
$ws.Range("E38").Value = "  -2.65%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5509"
$ws.Range("E39").Value = "  -3.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.407"
$ws.Range("E40").Value = "  -3.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1756"
$ws.Range("E41").Value = "  -2.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.851"
$ws.Range("E42").Value = "  +3.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.318"
$ws.Range("E43").Value = "  -3.24%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5170"
$ws.Range("E44").Value = "  -2.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.27"
$ws.Range("E45").Value = "  -3.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06834"
$ws.Range("E46").Value = "  -1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000002604"
$ws.Range("E47").Value = "  -6.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.041"
$ws.Range("E48").Value = "  -5.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.88"
$ws.Range("E49").Value = "  -2.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.779"
$ws.Range("E50").Value = "  -2.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.001"
$ws.Range("E51").Value = "  -0.21%  "
